$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level view tweak ---
$excel.ActiveWindow.TabRatio = 8

# --- Header row: regroup "Resting Rate"/"Start ramp"/"End ramp" into 3 groups ---
# (fall / rise / plateau). Column C becomes "Number of trials", column D becomes
# "Resting Rate", and two more (Resting Rate, Start ramp, End ramp) groups are
# appended after the existing one, shifting the data right by two columns.
$ws.Range("C1").Value = "Number of trials"
$ws.Range("D1").Value = "Resting Rate"
$ws.Range("G1").Value = "Resting Rate"
$ws.Range("H1").Value = "Start ramp"
$ws.Range("I1").Value = "End ramp"
$ws.Range("J1").Value = "Resting Rate"
$ws.Range("K1").Value = "Start ramp"
$ws.Range("L1").Value = "End ramp"

# --- Data rows: write the full new grid directly (H..L values taken/derived from the old G..J columns) ---
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 119400
$ws.Range("I2").Value = 121200
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 199600
$ws.Range("L2").Value = 201400

$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 2976
$ws.Range("I3").Value = 5271
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 489400
$ws.Range("L3").Value = 491500

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 192500
$ws.Range("I4").Value = 195400
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 315800
$ws.Range("L4").Value = 318500

$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 109400
$ws.Range("I5").Value = 111200
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 189600
$ws.Range("L5").Value = 191300

$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 125400
$ws.Range("I6").Value = 127400
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 215100
$ws.Range("L6").Value = 216800

$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 87050
$ws.Range("I7").Value = 89250

# --- Remove the old leftover "old"/stdev helper block in rows 22-28 ---
$ws.Range("A22:A28").EntireRow.Delete()

# --- Column widths: C widened, D/E unified to the same (old E) width ---
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 19.333333333333332
$ws.Columns.Item(5).ColumnWidth = 19.333333333333332

# --- Selection moves to L15 ---
$ws.Range("L15").Select()

Write-Output "done"
